# This script applies a set of individual cell value updates to the
# "Betfair Back/Lay" odds worksheet, as captured by the source diff.
# Each line updates exactly one cell identified by its A1-style reference
# on the active worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("P2").Value = 1.88
$ws.Range("Q4").Value = 1.96
$ws.Range("G5").Value = 1.63
$ws.Range("N5").Value = 5.7
$ws.Range("P5").Value = 2.58
$ws.Range("Q5").Value = 1.6
$ws.Range("R5").Value = 1.64
$ws.Range("T5").Value = 1.68
$ws.Range("U5").Value = 2.38
$ws.Range("AC5").Value = 10
$ws.Range("AK5").Value = 15.5
$ws.Range("AN5").Value = 6.6
$ws.Range("I6").Value = 2.3
$ws.Range("Y6").Value = 12
$ws.Range("AK6").Value = 38
$ws.Range("F7").Value = 5.4
$ws.Range("H7").Value = 1.71
$ws.Range("J7").Value = 4.3
$ws.Range("T7").Value = 1.84
$ws.Range("AE7").Value = 17
$ws.Range("AN7").Value = 1000
$ws.Range("F8").Value = 2.18
$ws.Range("H8").Value = 3.9
$ws.Range("I8").Value = 4
$ws.Range("J8").Value = 3.4
$ws.Range("T8").Value = 1.99
$ws.Range("X8").Value = 11.5
$ws.Range("AC8").Value = 7.6
$ws.Range("AE8").Value = 1000
$ws.Range("AJ8").Value = 30
$ws.Range("AL8").Value = 48
$ws.Range("AM8").Value = 130
$ws.Range("J9").Value = 3.3
$ws.Range("F10").Value = 1.58
$ws.Range("H10").Value = 6.2
$ws.Range("Q10").Value = 1.94
$ws.Range("F11").Value = 2.88
$ws.Range("H11").Value = 2.6
$ws.Range("I11").Value = 2.86
$ws.Range("J11").Value = 3.2
$ws.Range("K11").Value = 3.5
$ws.Range("P11").Value = 1.74
$ws.Range("Q11").Value = 2.12
$ws.Range("F12").Value = 1.94
$ws.Range("G12").Value = 2.06
$ws.Range("J12").Value = 3.7
$ws.Range("K12").Value = 4
$ws.Range("F13").Value = 1.44
$ws.Range("H13").Value = 6.6
$ws.Range("Q13").Value = 1.4
$ws.Range("P14").Value = 2.42
$ws.Range("Q14").Value = 1.42
$ws.Range("Q15").Value = 1.56
$ws.Range("F16").Value = 2.52
$ws.Range("H16").Value = 2.68
$ws.Range("I16").Value = 3.45
$ws.Range("J16").Value = 3.05
$ws.Range("K16").Value = 4.4
$ws.Range("P16").Value = 1.93
$ws.Range("Q16").Value = 1.86
$ws.Range("I17").Value = 1.66
$ws.Range("J17").Value = 4.1
$ws.Range("P17").Value = 2.5
$ws.Range("F18").Value = 1.62
$ws.Range("H18").Value = 5.2
$ws.Range("K18").Value = 4.7
$ws.Range("Q19").Value = 1.47
$ws.Range("G21").Value = 1.6
$ws.Range("P21").Value = 2.2
$ws.Range("F24").Value = 1.98
$ws.Range("H24").Value = 3.15
$ws.Range("H25").Value = 7.6
$ws.Range("U25").Value = 1.8
$ws.Range("AH25").Value = 28
$ws.Range("AI25").Value = 150
$ws.Range("AO25").Value = 260
$ws.Range("G26").Value = 1.39
$ws.Range("J26").Value = 5.7
$ws.Range("P26").Value = 2.5
$ws.Range("X26").Value = 25
$ws.Range("AH26").Value = 26
$ws.Range("AN26").Value = 5.2
$ws.Range("F27").Value = 2.02
$ws.Range("H27").Value = 3.1
$ws.Range("J27").Value = 3
$ws.Range("N27").Value = 1.67
$ws.Range("P27").Value = 1.67
$ws.Range("Q27").Value = 1.91
$ws.Range("S27").Value = 3.25
$ws.Range("V27").Value = 1.28
$ws.Range("W27").Value = 1.61
